# Tolerancia_CiruelaCan.xlsx update
# - Rename shared product-line label "CANDY PLUMS" -> "CANDY" (used by C2:C10)
# - Move the active selection from the header row (A1:XFD1) to cell C5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every row (2-10) in column C references the same shared string "CANDY PLUMS".
# Updating the whole range rewrites that shared string in place.
$ws.Range("C2:C10").Value = "CANDY"

# Move/collapse the saved selection to C5 (single cell, not the whole header row).
$ws.Range("C5").Select() | Out-Null
